# Restored from revision #c308e71254e0d1b2efda56730749fc766a931d86.TEST
# Author: admin. Type: SAVE.
#
# Update the "Integer min" rule value for R20 (row 10, column C) from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1

